$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.664.31'
$ws.Range('E2').Value = '  +0.83%  '
$ws.Range('D3').Value = '1.826.65'
$ws.Range('E3').Value = '  +1.41%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.009'
$ws.Range('E4').Value = '  +0.37%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.56'
$ws.Range('E6').Value = '  +0.61%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4700'
$ws.Range('E7').Value = '  +3.40%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3601'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07144'
$ws.Range('E9').Value = '  +0.60%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9352'
$ws.Range('E10').Value = '  +5.28%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07667'
$ws.Range('E11').Value = '  -1.40%  '
$ws.Range('E12').Value = '  -0.31%  '
$ws.Range('D13').Value = '1.865.14'
$ws.Range('E13').Value = '  +2.96%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.259'
$ws.Range('E14').Value = '  -0.70%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.356'
$ws.Range('E15').Value = '  +0.27%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '87.57'
$ws.Range('E16').Value = '  +2.43%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.010'
$ws.Range('E17').Value = '  +0.41%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008573'
$ws.Range('E18').Value = '  +0.62%  '
$ws.Range('E19').Value = '  +0.36%  '
$ws.Range('D20').Value = '26.672.03'
$ws.Range('E20').Value = '  +0.84%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.25'
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.031'
$ws.Range('E22').Value = '  +1.14%  '
$ws.Range('D23').Value = '2.064.16'
$ws.Range('E23').Value = '  +1.28%  '
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.919'
$ws.Range('E25').Value = '  -2.44%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '152.42'
$ws.Range('E26').Value = '  +0.95%  '
$ws.Range('E27').Value = '  +0.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.996'
$ws.Range('E28').Value = '  -2.16%  '
$ws.Range('E29').Value = '  +1.32%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.876'
$ws.Range('E30').Value = '  +0.23%  '
$ws.Range('E31').Value = '  +1.48%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.174'
$ws.Range('E32').Value = '  +1.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.836'
$ws.Range('E33').Value = '  +0.68%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.163'
$ws.Range('E34').Value = '  +4.82%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7364'
$ws.Range('E35').Value = '  +1.68%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.434'
$ws.Range('E36').Value = '  -0.48%  '
$ws.Range('E37').Value = '  +0.70%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01930'
$ws.Range('E38').Value = '  -0.38%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.936'
$ws.Range('E39').Value = '  +2.57%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.05157'
$ws.Range('E40').Value = '  +1.41%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.862'
$ws.Range('E41').Value = '  +0.60%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5065'
$ws.Range('E42').Value = '  -0.91%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1498'
$ws.Range('E43').Value = '  -1.33%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.092'
$ws.Range('E44').Value = '  +0.38%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.009'
$ws.Range('E45').Value = '  +0.51%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4648'
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.02'
$ws.Range('E47').Value = '  +0.84%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '98.57'
$ws.Range('E48').Value = '  -2.46%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.577'
$ws.Range('E49').Value = '  +0.42%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06036'
$ws.Range('E50').Value = '  +0.88%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '63.82'
$ws.Range('E51').Value = '  -0.01%  '
